$wb = $excel.ActiveWorkbook

# --- Sheet2 (TestData): fix "RunMode" -> "Runmode" typo on the two section headers ---
$wsData = $wb.Worksheets.Item("TestData")
$wsData.Range("A2").Value = "Runmode"
$wsData.Range("A8").Value = "Runmode"

# --- Add a new data row (11) duplicating the Amit Jena / Dollar row (10) ---
$wsData.Range("A11").Value = "Y"
$wsData.Range("B11").Value = "Amit Jena"
$wsData.Range("C11").Value = "Dollar"

# --- View state: TestData becomes the active sheet/tab, with A11:C11 selected ---
$wsData.Activate()
$excel.ActiveWindow.ScrollRow = 6
$excel.ActiveWindow.ScrollColumn = 1
$wsData.Range("A11:C11").Select()
